# Auto-generated edit script: update crypto price/volume table per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.660.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.791.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "432.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "

# Row 10
$ws.Range("E10").Value = "  -10.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.402.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.97%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.46%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.836.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("E18").Value = "  +2.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.732.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "410.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "

# Row 26
$ws.Range("E26").Value = "  +5.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +33.90%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.77%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "722.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.08%  "

# Row 31
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.56%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.136"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.56%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.46%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +28.98%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0474"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.69%  "

# Row 40
$ws.Range("E40").Value = "  +39.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0688"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.66%  "

# Row 43
$ws.Range("E43").Value = "  +3.36%  "

# Row 44
$ws.Range("E44").Value = "  +0.37%  "

# Row 45
$ws.Range("E45").Value = "  +2.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.322"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.42%  "

# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.57%  "

# Row 48
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("E49").Value = "  -0.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.32%  "

# Row 51
$ws.Range("E51").Value = "  -1.02%  "
